$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Entiat River Potato 07" row (original row 5).
# This shifts all subsequent rows up by one.
$ws.Rows.Item(5).Delete()

# Delete the "Methow River Thompson 09" row (original row 10, now row 9
# after the first deletion). This shifts remaining rows up by one more.
$ws.Rows.Item(9).Delete()

# After the two deletions, rows now look like:
#  5  Entiat River Potato 08
#  6  Methow River Fawn 04
#  7  Methow River Thompson 07
#  8  Methow River Thompson 08
#  9  Nason Creek Lower 01
# 10  Nason Creek Lower 02
# 11  Nason Creek Lower 03
# 12  Twisp River Middle 01
# 13  Twisp River Middle 02
# 14  Twisp River Middle 06

# Apply the remaining score updates to match the refreshed data.

# Row 5: Entiat River Potato 08
$ws.Range("K5").Value = 5
$ws.Range("P5").Value = 3
$ws.Range("R5").Value = 3
$ws.Range("T5").Value = 33
$ws.Range("U5").Value = 0.7333333333333333

# Row 7: Methow River Thompson 07
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 1
$ws.Range("T7").Value = 32
$ws.Range("U7").Value = 0.7111111111111111

# Row 8: Methow River Thompson 08
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("T8").Value = 32
$ws.Range("U8").Value = 0.7111111111111111

# Row 9: Nason Creek Lower 01
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = 3
$ws.Range("T9").Value = 36
$ws.Range("U9").Value = 0.8

# Row 11: Nason Creek Lower 03
$ws.Range("P11").Value = 3
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = 3
$ws.Range("T11").Value = 34
$ws.Range("U11").Value = 0.7555555555555555
$ws.Range("V11").Value = 5
